$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 370
$ws.Range("J2").Value = 1559
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 429
$ws.Range("M2").Value = 37
$ws.Range("N2").Value = 286
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 24
$ws.Range("S2").Value = 178
$ws.Range("T2").Value = 279
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 2576
$ws.Range("X2").Value = 2473
$ws.Range("Z2").Value = 46
$ws.Range("AA2").Value = 11
